# Add the word "option" for clarity in the wireless command sentence:
#   "...Recommend command for wireless from the devs..."
#   becomes
#   "...Recommend command for wireless option from the devs..."
$d = $word.ActiveDocument

$full = $d.Content.Text
$marker = "wireless from the devs"
$idx = $full.IndexOf($marker)

if ($idx -ge 0) {
    # Insert the single word precisely between "wireless " and "from the devs"
    # using a zero-length Range so we don't disturb any other text/formatting.
    $insertPos = $idx + "wireless ".Length
    $r = $d.Range($insertPos, $insertPos)
    $r.InsertBefore("option ")
    Write-Host "Inserted 'option' at position $insertPos"
} else {
    # Fallback: use Find/Replace in case the exact offsets could not be located.
    $find = $d.Content.Find
    $find.ClearFormatting()
    $found = $find.Execute("wireless from the devs", $true, $false, $false, $false, $false, $true, 1, $false, "wireless option from the devs", 2)
    Write-Host "Fallback Find/Replace executed: $found"
}
